# Aufgabe in Teil 1-6 zerlegt
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Helper: within a TextRange that may contain several paragraphs, split
# the FIRST paragraph's single run of text "prefix+suffix" into two
# separate runs "prefix" / "suffix" (same rPr), without touching any
# other paragraph in the text frame.
#   $fullTextRange : TextFrame.TextRange of the shape
#   $prefix        : text that stays in run 1
#   $suffix        : text that becomes run 2
# The first paragraph's current text must equal prefix+suffix already
# (or will be set to it).
# ---------------------------------------------------------------------
function Split-FirstParagraphRun($fullTextRange, $prefix, $suffix) {
    $combined = $prefix + $suffix
    $whole = $fullTextRange.Characters(1, $combined.Length)
    $whole.Text = $combined
    $part1 = $fullTextRange.Characters(1, $prefix.Length)
    $part1.Text = $prefix
}

# 1) Slide 14: "Exemplare – Übung 1" -> "Exemplare – Übung " + "(Teil 1)"
$s14 = $p.Slides.Item(14)
$title14 = $s14.Shapes.Title
Split-FirstParagraphRun $title14.TextFrame.TextRange "Exemplare – Übung " "(Teil 1)"

# 2) Slide 2: "Objektorientierte Programmierung" -> "Objektorientierte " + "Programmierung"
$s2 = $p.Slides.Item(2)
$subtitle2 = $s2.Shapes.Item(5)
Split-FirstParagraphRun $subtitle2.TextFrame.TextRange "Objektorientierte " "Programmierung"

# 2b) Slide 2: add new "Rechteck 1" shape with GitHub hyperlink
$rect = $s2.Shapes.AddShape(1, 147.39842519685038, 383.4668503937008, 515.9628346456693, 27.556535433070867)
$rect.Name = "Rechteck 1"
$rect.TextFrame.WordWrap = -1
$rect.TextFrame.AutoSize = 1

$url = "https://github.com/nordakademie-einfuehrung-java/uebung_6"
$rectTr = $rect.TextFrame.TextRange
$rectTr.Text = "https://"
$rectPart2 = $rectTr.InsertAfter("github.com/nordakademie-einfuehrung-java/uebung_6")

$rectTr.ActionSettings(1).Hyperlink.Address = $url
$rectPart2.ActionSettings(1).Hyperlink.Address = $url

# Restore exact target height (AutoSize recalculated height on text change)
$rect.Height = 27.556535433070867

# 3) Slide 20: "Exemplare – Übung 2" -> "Exemplare – Übung " + "(Teil 2)"
$s20 = $p.Slides.Item(20)
$title20 = $s20.Shapes.Title
Split-FirstParagraphRun $title20.TextFrame.TextRange "Exemplare – Übung " "(Teil 2)"

# 4) Slide 20: merge "Implementieren Sie folgende Methode, um den Tank " + "aufzufüllen: "
#    into a single run.
$body20 = $s20.Shapes.Item(5)
$tr20 = $body20.TextFrame.TextRange
$full20 = $tr20.Text
$search20 = "Implementieren Sie folgende Methode, um den Tank "
$tail20 = "aufzufüllen: "
$idx20 = $full20.IndexOf($search20)
$seg20 = $tr20.Characters($idx20 + 1, $search20.Length + $tail20.Length)
$seg20.Text = $search20 + $tail20

# 5) Slide 26: "Übung 1" -> "Übung " + "(Teil 3)"
$s26 = $p.Slides.Item(26)
$title26 = $s26.Shapes.Title
Split-FirstParagraphRun $title26.TextFrame.TextRange "Übung " "(Teil 3)"

# 6) Slide 27: "Übung 2" -> "Übung " + "(Teil 4)"
$s27 = $p.Slides.Item(27)
$title27 = $s27.Shapes.Title
Split-FirstParagraphRun $title27.TextFrame.TextRange "Übung " "(Teil 4)"

# 7) Slide 28: "Übung 3" -> "Übung " + "(Teil 5)"
$s28 = $p.Slides.Item(28)
$title28 = $s28.Shapes.Title
Split-FirstParagraphRun $title28.TextFrame.TextRange "Übung " "(Teil 5)"

# 8) Slide 29: "Übung 4" -> "Übung " + "(Teil 6)"
$s29 = $p.Slides.Item(29)
$title29 = $s29.Shapes.Title
Split-FirstParagraphRun $title29.TextFrame.TextRange "Übung " "(Teil 6)"

Write-Output "edit.ps1 applied"
